# "Use locale code instead": the Locale column value changes from the
# English language name to its locale code, and the two header cells
# (Slug, Locale) plus the row-2 "Slug" cell pick up the plain body-text
# formatting already used elsewhere in the sheet (e.g. column G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2, column I ("Locale") - English -> en
$ws.Range("I2").Value = "en"

# Copy the already-normalized formatting from column G onto the header
# cells H1:I1 and the row-2 "Slug" cell H2.
$ws.Range("G1").Copy()
$ws.Range("H1:I1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("G2").Copy()
$ws.Range("H2").PasteSpecial(-4122)      # xlPasteFormats

$excel.CutCopyMode = 0
